# Apply updated cryptocurrency market data (price + 1h volume change) to sheet1.
# Source rows are keyed by their original spreadsheet row number; a handful of rows
# also have their Coin name / Link swapped with a neighboring row (re-ranking).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.811.48"
$ws.Range("E2").Value = "  -0.39%  "

# Row 3
$ws.Range("D3").Value = "1.866.20"
$ws.Range("E3").Value = "  -1.42%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D5").Value = "0.7333"
$ws.Range("E5").Value = "  -5.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D6").Value = "241.52"
$ws.Range("E6").Value = "  -0.98%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.22%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D8").Value = "0.3095"
$ws.Range("E8").Value = "  -1.18%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D9").Value = "24.55"
$ws.Range("E9").Value = "  -4.06%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D10").Value = "0.07056"
$ws.Range("E10").Value = "  -3.75%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D11").Value = "0.08390"
$ws.Range("E11").Value = "  +4.17%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.898.83"
$ws.Range("E12").Value = "  +4.43%  "

# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D13").Value = "0.7465"
$ws.Range("E13").Value = "  -3.13%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D14").Value = "5.318"
$ws.Range("E14").Value = "  -3.11%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D15").Value = "92.13"
$ws.Range("E15").Value = "  -2.04%  "

# Row 16
$ws.Range("D16").Value = "29.816.01"
$ws.Range("E16").Value = "  -0.19%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D17").Value = "5.988"
$ws.Range("E17").Value = "  -3.56%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D18").Value = "13.52"
$ws.Range("E18").Value = "  -3.37%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D19").Value = "0.000007782"
$ws.Range("E19").Value = "  -0.82%  "

# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D20").Value = "239.49"
$ws.Range("E20").Value = "  -2.81%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  +0.29%  "

# Row 22
$ws.Range("D22").Value = "2.134.03"
$ws.Range("E22").Value = "  +2.52%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D24").Value = "7.881"
$ws.Range("E24").Value = "  -3.16%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D25").Value = "0.1560"
$ws.Range("E25").Value = "  -0.82%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D26").Value = "9.267"
$ws.Range("E26").Value = "  -1.74%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D27").Value = "162.45"
$ws.Range("E27").Value = "  +0.24%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D28").Value = "18.50"
$ws.Range("E28").Value = "  -1.34%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D29").Value = "1.995"
$ws.Range("E29").Value = "  -1.31%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D30").Value = "1.494"
$ws.Range("E30").Value = "  +5.01%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D31").Value = "1.526"
$ws.Range("E31").Value = "  -0.90%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D32").Value = "4.451"
$ws.Range("E32").Value = "  -0.53%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D33").Value = "4.123"
$ws.Range("E33").Value = "  +1.35%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D34").Value = "0.05358"
$ws.Range("E34").Value = "  -3.37%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D35").Value = "1.227"
$ws.Range("E35").Value = "  -0.67%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D36").Value = "0.7431"
$ws.Range("E36").Value = "  -0.82%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D37").Value = "0.9990"
$ws.Range("E37").Value = "  -0.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D38").Value = "2.699"
$ws.Range("E38").Value = "  +0.68%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D39").Value = "0.01928"
$ws.Range("E39").Value = "  +0.13%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D40").Value = "2.768"
$ws.Range("E40").Value = "  -0.74%  "

# Row 41
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.108.19"
$ws.Range("E41").Value = "  +1.38%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D42").Value = "0.4426"
$ws.Range("E42").Value = "  -1.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D43").Value = "5.995"
$ws.Range("E43").Value = "  -0.18%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D44").Value = "71.89"
$ws.Range("E44").Value = "  -3.15%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D45").Value = "0.8638"
$ws.Range("E45").Value = "  +1.61%  "

# Row 46
$ws.Range("E46").Value = "  +0.23%  "

# Row 47
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D47").Value = "102.03"
$ws.Range("E47").Value = "  -0.39%  "

# Row 48
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D48").Value = "7.670"
$ws.Range("E48").Value = "  +1.67%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D49").Value = "1.830"
$ws.Range("E49").Value = "  -2.91%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"  # keep exact text (avoid numeric auto-conversion)
$ws.Range("D50").Value = "2.998"
$ws.Range("E50").Value = "  +0.05%  "

# Row 51
$ws.Range("D51").Value = "2.033.30"
$ws.Range("E51").Value = "  +1.19%  "
